$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skill_auto")

# Fix SKILL503 detail text: "ターン終了時にHPが10%回復" -> "ターン終了時にHPが10回復"
$ws.Range("O5").Value = "ターン終了時にHPが10回復"

# Add new row 8 for SKILL506 "精神統一" (AUTO, TP +5 recovery at end of turn).
# Copy formatting from the row above (row 7) so borders/fonts/fills match the
# rest of the table, then fill in the row height and the actual cell values.
$ws.Range("A7:O7").Copy()
$ws.Range("A8:O8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Rows.Item(8).RowHeight = 20

$ws.Range("A8").Value = "SKILL506"
$ws.Range("B8").Value = "精神統一"
$ws.Range("C8").Value = "AUTO"
$ws.Range("G8").Value = 5
$ws.Range("O8").Value = "ターン終了時にTPが5回復"
